# output stage with 1uf real caps
# Insert a new worksheet (becomes "Sheet6" / third tab) right after "Sheet5"
# and before "Sheet2", populated with real-capacitor data for the output
# stage (MKP1840 / PHE426 / R75 / MKP4 / MKP10 parts).

$wb = $excel.ActiveWorkbook

$sheet5 = $wb.Worksheets.Item("Sheet5")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet5)
$ws.Name = "Sheet6"

# Column widths (closest achievable values to the authored 13.6640625 / 5)
$ws.Columns.Item(1).ColumnWidth = 12.83
$ws.Columns.Item(2).ColumnWidth = 4.17

# Header row
$ws.Range("A1").Value = "max V @ 40kHz"
$ws.Range("C1").Value = "160VDC"
$ws.Range("D1").Value = "250VDC"
$ws.Range("E1").Value = "400VDC"
$ws.Range("F1").Value = "630VDC"

# MKP1840 block
$ws.Range("A2").Value = "MKP1840"
$ws.Range("B2").Value = 0.47
$ws.Range("C2").Value = 33
$ws.Range("D2").Value = 33
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 43

$ws.Range("B3").Value = 0.68

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 22
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 30

$ws.Range("B5").Value = 1.5

# PHE426 block
$ws.Range("A7").Value = "PHE426"
$ws.Range("B7").Value = 0.47

$ws.Range("B8").Value = 0.56

$ws.Range("B9").Value = 0.68

$ws.Range("B10").Value = 0.82
$ws.Range("E10").Value = 40

$ws.Range("B11").Value = 1
$ws.Range("D11").Formula = "=21*SQRT(2.2)"
$ws.Range("E11").Formula = "=E10/SQRT(B11/B10)"
$ws.Range("F11").Value = 32

$ws.Range("B12").Value = 1.2

$ws.Range("B13").Value = 1.5

# R75 block
$ws.Range("A15").Value = "R75"
$ws.Range("B15").Value = 0.47
$ws.Range("E15").Value = 45

$ws.Range("B16").Value = 0.56

$ws.Range("B17").Value = 0.68
$ws.Range("D17").Value = 33

$ws.Range("B18").Value = 0.82

$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 29
$ws.Range("E19").Formula = "=E15/SQRT(B19/B15)"
$ws.Range("F19").Value = 30

$ws.Range("B20").Value = 1.2

$ws.Range("B21").Value = 1.5
$ws.Range("D21").Value = 21

# MKP4 block
$ws.Range("A23").Value = "MKP4"
$ws.Range("B23").Value = 0.47

$ws.Range("B24").Value = 0.68

$ws.Range("B25").Value = 1
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 14

$ws.Range("B26").Value = 1.5

# MKP10 block
$ws.Range("A28").Value = "MKP10"
$ws.Range("B28").Value = 0.47
$ws.Range("D28").Value = 23
$ws.Range("E28").Value = 30
$ws.Range("F28").Value = 60

$ws.Range("B29").Value = 0.68

$ws.Range("B30").Value = 1
$ws.Range("D30").Value = 16
$ws.Range("E30").Value = 22
$ws.Range("F30").Value = 32

$ws.Range("B31").Value = 1.5

# Match the authored selection/active cell on the new sheet
$ws.Range("D11").Select()
